$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, derived from the commit diff.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the workbook's inlineStr cells) instead of
# auto-converting numeric-/percent-looking strings into numbers.
# The Style is reset to "Normal" afterwards so the forced quote-prefix
# text formatting does not linger on the cell.
$updates = [ordered]@{
    "D2" = "292.46"
    "E2" = "0.00%"
    "D3" = "40.46"
    "E3" = "0.28%"
    "D4" = "5.014"
    "E4" = "-0.68%"
    "D5" = "0.07290"
    "E5" = "-1.59%"
    "B6" = "GateToken"
    "C6" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D6" = "4.282"
    "E6" = "-0.93%"
    "B7" = "FTXToken"
    "C7" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D7" = "1.559"
    "E7" = "-1.35%"
    "B8" = "MXToken"
    "C8" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D8" = "0.9271"
    "E8" = "0.26%"
    "B9" = "BTSEToken"
    "C9" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D9" = "2.333"
    "E9" = "-3.60%"
    "B10" = "LiechtensteinCryptoassetsExchange"
    "C10" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D10" = "0.1158"
    "E10" = "-1.48%"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D11" = "0.1763"
    "E11" = "0.87%"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D12" = "0.04367"
    "E12" = "4.31%"
    "B13" = "MandalaExchangeToken"
    "C13" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D13" = "0.08724"
    "E13" = "0.01%"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D14" = "0.1055"
    "E14" = "0.12%"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D15" = "0.001277"
    "E15" = "0.78%"
    "D16" = "0.005973"
    "E16" = "1.19%"
    "E17" = "-0.67%"
    "B18" = "BitpandaEcosystemToken"
    "C18" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "D18" = "0.3288"
    "E18" = "-1.79%"
    "B19" = "MCDex"
    "C19" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D19" = "7.908"
    "E19" = "3.33%"
    "B20" = "ProBitToken"
    "C20" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "D20" = "0.1391"
    "E20" = "2.26%"
    "B21" = "ZBToken"
    "C21" = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
    "D21" = "0.2774"
    "E21" = "-1.72%"
    "B22" = "CoinExToken"
    "C22" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "D22" = "0.03924"
    "E22" = "1.42%"
    "D23" = "0.001262"
    "D24" = "0.003668"
    "E24" = "4.86%"
    "E25" = "-8.34%"
    "D26" = "0.0003727"
    "E26" = "-0.64%"
    "D38" = "0.02299"
    "E38" = "-0.56%"
    "D39" = "0.05070"
    "E39" = "1.46%"
    "D40" = "0.005723"
    "E40" = "37.28%"
    "D41" = "0.007854"
    "E41" = "1.40%"
    "D42" = "0.1285"
    "E42" = "0.55%"
    "D43" = "0.007385"
    "E43" = "-0.85%"
    "D44" = "0.007280"
    "E44" = "1.91%"
    "D45" = "0.2902"
    "E45" = "-8.98%"
    "D46" = "0.00006208"
    "E46" = "-7.28%"
    "E47" = "-0.64%"
    "D48" = "0.04856"
    "E48" = "-80.72%"
    "D49" = "0.00002103"
    "E49" = "-0.64%"
    "D50" = "0.0002003"
    "E50" = "-0.64%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.Value = "'" + $updates[$cell]
    $range.Style = "Normal"
}
